{"js": "const replacements = [\n  ['2026-01-03 Saturday', '2026-01-04 Sunday'],\n  ['839\u00f79=93, 2', '682\u00f74=170, 2'],\n  ['963\u00f74=240, 3', '660\u00f78=82, 4'],\n  ['907\u00f76=151, 1', '325\u00f77=46, 3'],\n  ['177\u00f75=35, 2', '700\u00f72=350, 0'],\n  ['302\u00f74=75, 2', '260\u00f73=86, 2'],\n  ['297\u00f75=59, 2', '995\u00f75=199, 0'],\n  ['668\u00f79=74, 2', '433\u00f76=72, 1'],\n  ['250\u00f78=31, 2', '321\u00f75=64, 1'],\n  ['422\u00f72=211, 0', '511\u00f75=102, 1'],\n  ['564\u00f73=188, 0', '872\u00f74=218, 0'],\n  ['927\u00f75=185, 2', '973\u00f77=139, 0'],\n  ['294\u00f76=49, 0', '844\u00f78=105, 4'],\n  ['190\u00f77=27, 1', '362\u00f73=120, 2'],\n  ['557\u00f77=79, 4', '817\u00f72=408, 1'],\n  ['509\u00f74=127, 1', '815\u00f78=101, 7'],\n  ['334\u00f72=167, 0', '396\u00f72=198, 0'],\n  ['642\u00f78=80, 2', '117\u00f78=14, 5'],\n  ['345\u00f79=38, 3', '394\u00f72=197, 0'],\n  ['562\u00f73=187, 1', '838\u00f76=139, 4'],\n  ['520\u00f78=65, 0', '271\u00f79=30, 1'],\n  ['813\u00f76=135, 3', '113\u00f75=22, 3'],\n  ['119\u00f72=59, 1', '895\u00f76=149, 1'],\n  ['694\u00f77=99, 1', '428\u00f72=214, 0'],\n  ['198\u00f74=49, 2', '629\u00f78=78, 5'],\n  ['168\u00f73=56, 0', '228\u00f74=57, 0'],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load('items');\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  para.load('text');\n}\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  const text = para.text;\n  for (const [before, after] of replacements) {\n    if (text === before) {\n      para.insertText(after, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "$replacements = @(\n    @('2026-01-03 Saturday', '2026-01-04 Sunday'),\n    @('839\u00f79=93, 2', '682\u00f74=170, 2'),\n    @('963\u00f74=240, 3', '660\u00f78=82, 4'),\n    @('907\u00f76=151, 1', '325\u00f77=46, 3'),\n    @('177\u00f75=35, 2', '700\u00f72=350, 0'),\n    @('302\u00f74=75, 2', '260\u00f73=86, 2'),\n    @('297\u00f75=59, 2', '995\u00f75=199, 0'),\n    @('668\u00f79=74, 2', '433\u00f76=72, 1'),\n    @('250\u00f78=31, 2', '321\u00f75=64, 1'),\n    @('422\u00f72=211, 0', '511\u00f75=102, 1'),\n    @('564\u00f73=188, 0', '872\u00f74=218, 0'),\n    @('927\u00f75=185, 2', '973\u00f77=139, 0'),\n    @('294\u00f76=49, 0', '844\u00f78=105, 4'),\n    @('190\u00f77=27, 1', '362\u00f73=120, 2'),\n    @('557\u00f77=79, 4', '817\u00f72=408, 1'),\n    @('509\u00f74=127, 1', '815\u00f78=101, 7'),\n    @('334\u00f72=167, 0', '396\u00f72=198, 0'),\n    @('642\u00f78=80, 2', '117\u00f78=14, 5'),\n    @('345\u00f79=38, 3', '394\u00f72=197, 0'),\n    @('562\u00f73=187, 1', '838\u00f76=139, 4'),\n    @('520\u00f78=65, 0', '271\u00f79=30, 1'),\n    @('813\u00f76=135, 3', '113\u00f75=22, 3'),\n    @('119\u00f72=59, 1', '895\u00f76=149, 1'),\n    @('694\u00f77=99, 1', '428\u00f72=214, 0'),\n    @('198\u00f74=49, 2', '629\u00f78=78, 5'),\n    @('168\u00f73=56, 0', '228\u00f74=57, 0')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $find.Execute(\n        $findText,      # FindText\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
